$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(10002, 10021, 3000121),
    @(10003, 10022, 3000122),
    @(10004, 10023, 3000123),
    @(10005, 10024, 3000124),
    @(10006, 10025, 3000125),
    @(10007, 10026, 3000126),
    @(10008, 10027, 3000127),
    @(10009, 10028, 3000128),
    @(10010, 10029, 3000129),
    @(10002, 10021, 3000130),
    @(10003, 10022, 3000131),
    @(10004, 10023, 3000132),
    @(10005, 10024, 3000133),
    @(10006, 10025, 3000134),
    @(10007, 10026, 3000135),
    @(10008, 10027, 3000136),
    @(10009, 10028, 3000137),
    @(10010, 10029, 3000138),
    @(10002, 10021, 3000139),
    @(10003, 10022, 3000140),
    @(10004, 10023, 3000141),
    @(10005, 10024, 3000142),
    @(10006, 10025, 3000143),
    @(10007, 10026, 3000144),
    @(10008, 10027, 3000145),
    @(10009, 10028, 3000146),
    @(10010, 10029, 3000147),
    @(10002, 10021, 3000148),
    @(10003, 10022, 3000149),
    @(10004, 10023, 3000150),
    @(10005, 10024, 3000151),
    @(10006, 10025, 3000152),
    @(10007, 10026, 3000153),
    @(10008, 10027, 3000154),
    @(10009, 10028, 3000155),
    @(10010, 10029, 3000156),
    @(10002, 10021, 3000157),
    @(10003, 10022, 3000158),
    @(10004, 10023, 3000159),
    @(10005, 10024, 3000160),
    @(10006, 10025, 3000161),
    @(10007, 10026, 3000162),
    @(10008, 10027, 3000163),
    @(10009, 10028, 3000164),
    @(10010, 10029, 3000165)
)

$r = 102
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin()"
    $ws.Cells.Item($r, 7).Value = "now()"
    $r++
}